$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Which minerals are on both the US Geological Survey and the Department of Energy's lists of critical minerals?"
$ws.Range("C4").Value = "Cobalt, Graphite, Lithium, Manganese and Nickel"

$ws.Range("A5").Value = "Which country is the leading producer of cobalt?"
$ws.Range("C5").Value = "The Democratic Republic of Congo"

$ws.Range("A6").Value = "Which country is the leading producer of Lithium?"
$ws.Range("C6").Value = "Australia"

$ws.Range("A7").Value = "Which country is the leading producer of rare earth minerals"
$ws.Range("C7").Value = "China"

$ws.Range("F12").Select()
